$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting from the row above (A15 uses the bold/border/centered style)
# so the new A16 cell picks up the same style, then set its value.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.632336121838838
$ws.Range("D16").Value = 1.039049773750829
$ws.Range("E16").Value = 0.8669497803913691
$ws.Range("F16").Value = 1.632336121838838
$ws.Range("G16").Value = 0.9407265842407487
$ws.Range("H16").Value = 1.56056169696832
$ws.Range("I16").Value = 0.9008325933958826
$ws.Range("J16").Value = 1.039049773750829
$ws.Range("K16").Value = 0.9529997770710992
$ws.Range("L16").Value = 1.292667949454969
$ws.Range("M16").Value = 1.156742758430998
